$d = $word.ActiveDocument
$style = $d.Styles("Default Paragraph Font")
$style.Visible = $true
